$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.238.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.707.78'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +2.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.707.19'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('E10').Value = '  +1.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.363'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.23'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.202.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.208.68'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.710.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.73'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '369.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.48'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.58%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.838.94'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '576.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.98'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.54%  '
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.60'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.83'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '157.98'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.378'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.81%  '
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0309'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.72'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.596'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '155.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('E51').Value = '  +3.37%  '
